# Generate Report for Handoff
# Refresh the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps for the
# files that were re-handed-off (status: "Handback transform failed" or
# "Ready for handoff") across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = Latest Handoff Date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-25 00:39:40"
}

# --- zh-cn sheet: column E = Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $zhCnRows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-25 00:39:34"
}

# --- de-de sheet: column E = Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $deDeRows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-25 00:39:40"
}
